# Inclusao da funcao de precos + logotipos do restaurante
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "preco" column header (D1), styled like the other headers (bold,
# bordered, centered) by copying the format of an existing header cell.
$ws.Range("D1").Value = "preco"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Updated quantities in column B
$ws.Range("B2").Value = 0
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 3000

# New "preco" values in column D for existing rows
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 44
$ws.Range("D6").Value = 200
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0

# New rows 10 and 11
$ws.Range("A10").Value = "tomate seco timy sache"
$ws.Range("B10").Value = 3000
$ws.Range("C10").Value = "g"
$ws.Range("D10").Value = 50

$ws.Range("A11").Value = "azeitona"
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = "Un"
$ws.Range("D11").Value = 12
